# "Add cantrals by cantons"
# The sheet used to have two header rows (row 1 and row 2) describing the
# plant measurements; this consolidates them into a single header row with
# new column labels (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1),
# (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year) and shifts the existing
# plant data rows up by one row to sit right under the new header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the second (old) header row; this shifts all the plant data rows
# up by one, right under row 1.
$ws.Rows("2:2").Delete()

# Wipe out whatever remains of the old row-1 header (values + styles) so we
# can build the new, single header row from scratch.
$ws.Rows("1:1").Clear()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Match the font used throughout the rest of the table for the numeric /
# unit columns of the new header.
$units = $ws.Range("F1:K1")
$units.Font.Name = "Arial"
$units.Font.Size = 9

# Match the selection left behind by the edit (the whole first data row).
$ws.Range("A2:K2").Select() | Out-Null
